$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new forecast row (row 20), reusing the date-cell formatting (borders,
# bold font, custom date number format) already used by column A in row 19.
$ws.Range("A19").Copy()
$ws.Range("A20").PasteSpecial(-4122)   # xlPasteFormats

$ws.Range("A20").Value = 45986
$ws.Range("B20").Value = 2025
$ws.Range("C20").Value = 1.049317648994741
$ws.Range("D20").Value = 2026
$ws.Range("E20").Value = 0.07146359800258573
